$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update localized/representation text values to the invariant/English forms
$ws.Range("F3").Value = "9/2/2010 12:00:00 AM"
$ws.Range("G3").Value = "2010-Sep-02"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "1,234.57"
$ws.Range("G5").NumberFormat = ""

# Widen column F from 18.7 to 21.7
# Note: the COM ColumnWidth setter here quantizes to a 1/6-character grid,
# so 20.83 (-> stored width ~21.67, closest reachable value to 21.7) is used.
$ws.Columns.Item(6).ColumnWidth = 20.83
